# --- Add a new "2022-Q1" worksheet (fund holdings), positioned right before "总计" ---
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Activate()
$q1 = $wb.Worksheets.Add()   # Add() inserts right before the active sheet, i.e. before "总计"
$q1.Name = "2022-Q1"

# Header row, matching the other quarterly sheets
$q1.Range("B1").Value = '基金代码'
$q1.Range("C1").Value = '基金名称'
$q1.Range("D1").Value = '基金规模'
$q1.Range("E1").Value = '股票总仓位'
$q1.Range("F1").Value = '仓位占比'
$q1.Range("G1").Value = '持有市值(亿元)'
$q1.Range("H1").Value = '仓位排名'
$q1.Range("B1:H1").Font.Bold = $true
$q1.Range("B1:H1").HorizontalAlignment = -4108   # xlCenter
$q1.Range("B1:H1").VerticalAlignment = -4160     # xlTop
$q1.Range("B1:H1").Borders.LineStyle = 1         # xlContinuous

# Data rows 2-13: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值/仓位排名
$q1.Range("A2").Value = 0
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = '005821'
$q1.Range("C2").Value = '万家新机遇龙头企业灵活配置混合'
$q1.Range("D2").Value = '23.23'
$q1.Range("E2").Value = '56.20'
$q1.Range("F2").Value = '2.92'
$q1.Range("G2").Value = '0.6783'
$q1.Range("H2").Value = 7

$q1.Range("A3").Value = 1
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = '013960'
$q1.Range("C3").Value = '万家新机遇成长一年持有期混合A'
$q1.Range("D3").Value = '13.29'
$q1.Range("E3").Value = '49.51'
$q1.Range("F3").Value = '3.05'
$q1.Range("G3").Value = '0.4053'
$q1.Range("H3").Value = 6

$q1.Range("A4").Value = 2
$q1.Range("B4:G4").NumberFormat = "@"
$q1.Range("B4").Value = '160642'
$q1.Range("C4").Value = '鹏华增瑞灵活配置混合(LOF)'
$q1.Range("D4").Value = '6.76'
$q1.Range("E4").Value = '91.34'
$q1.Range("F4").Value = '3.48'
$q1.Range("G4").Value = '0.2352'
$q1.Range("H4").Value = 9

$q1.Range("A5").Value = 3
$q1.Range("B5:G5").NumberFormat = "@"
$q1.Range("B5").Value = '168501'
$q1.Range("C5").Value = '北信瑞丰产业升级多策略混合'
$q1.Range("D5").Value = '4.42'
$q1.Range("E5").Value = '94.11'
$q1.Range("F5").Value = '3.75'
$q1.Range("G5").Value = '0.1658'
$q1.Range("H5").Value = 8

$q1.Range("A6").Value = 4
$q1.Range("B6:G6").NumberFormat = "@"
$q1.Range("B6").Value = '013961'
$q1.Range("C6").Value = '万家新机遇成长一年持有期混合C'
$q1.Range("D6").Value = '3.13'
$q1.Range("E6").Value = '49.51'
$q1.Range("F6").Value = '3.05'
$q1.Range("G6").Value = '0.0955'
$q1.Range("H6").Value = 6

$q1.Range("A7").Value = 5
$q1.Range("B7:G7").NumberFormat = "@"
$q1.Range("B7").Value = '005933'
$q1.Range("C7").Value = '新疆前海联合先进制造灵活配置混合A'
$q1.Range("D7").Value = '1.26'
$q1.Range("E7").Value = '89.79'
$q1.Range("F7").Value = '7.38'
$q1.Range("G7").Value = '0.0930'
$q1.Range("H7").Value = 2

$q1.Range("A8").Value = 6
$q1.Range("B8:G8").NumberFormat = "@"
$q1.Range("B8").Value = '001396'
$q1.Range("C8").Value = '建信互联网+产业升级股票'
$q1.Range("D8").Value = '2.50'
$q1.Range("E8").Value = '84.08'
$q1.Range("F8").Value = '3.22'
$q1.Range("G8").Value = '0.0805'
$q1.Range("H8").Value = 10

$q1.Range("A9").Value = 7
$q1.Range("B9:G9").NumberFormat = "@"
$q1.Range("B9").Value = '001056'
$q1.Range("C9").Value = '北信瑞丰健康生活主题灵活配置混合'
$q1.Range("D9").Value = '1.64'
$q1.Range("E9").Value = '86.03'
$q1.Range("F9").Value = '3.79'
$q1.Range("G9").Value = '0.0622'
$q1.Range("H9").Value = 7

$q1.Range("A10").Value = 8
$q1.Range("B10:G10").NumberFormat = "@"
$q1.Range("B10").Value = '004128'
$q1.Range("C10").Value = '新疆前海联合泳隆灵活配置混合A'
$q1.Range("D10").Value = '0.86'
$q1.Range("E10").Value = '91.05'
$q1.Range("F10").Value = '6.46'
$q1.Range("G10").Value = '0.0556'
$q1.Range("H10").Value = 1

$q1.Range("A11").Value = 9
$q1.Range("B11:G11").NumberFormat = "@"
$q1.Range("B11").Value = '007040'
$q1.Range("C11").Value = '新疆前海联合泳隆灵活配置混合C'
$q1.Range("D11").Value = '0.82'
$q1.Range("E11").Value = '91.05'
$q1.Range("F11").Value = '6.46'
$q1.Range("G11").Value = '0.0530'
$q1.Range("H11").Value = 1

$q1.Range("A12").Value = 10
$q1.Range("B12:G12").NumberFormat = "@"
$q1.Range("B12").Value = '005934'
$q1.Range("C12").Value = '新疆前海联合先进制造灵活配置混合C'
$q1.Range("D12").Value = '0.09'
$q1.Range("E12").Value = '89.79'
$q1.Range("F12").Value = '7.38'
$q1.Range("G12").Value = '0.0066'
$q1.Range("H12").Value = 2

$q1.Range("A13").Value = 11
$q1.Range("B13:G13").NumberFormat = "@"
$q1.Range("B13").Value = '001866'
$q1.Range("C13").Value = '北信瑞丰新成长灵活配置混合'
$q1.Range("D13").Value = '0.07'
$q1.Range("E13").Value = '94.21'
$q1.Range("F13").Value = '4.92'
$q1.Range("G13").Value = '0.0034'
$q1.Range("H13").Value = 4

# Column-A index cells use the same bold/bordered/centered look as the header
$q1.Range("A2:A13").Font.Bold = $true
$q1.Range("A2:A13").HorizontalAlignment = -4108
$q1.Range("A2:A13").VerticalAlignment = -4160
$q1.Range("A2:A13").Borders.LineStyle = 1
# --- Insert the new "2022-Q1" summary row at the top of the "总计" data table ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Undo any formatting the row-insert inherited from the header row above it -
# the new data row (like every other data row in this table) carries no
# explicit cell style except for the bold/bordered index cell in column A.
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 1.93

# Match the existing look of the other index cells in column A (bold/bordered/centered)
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

# Renumber the pre-existing rows' index column (A) to stay sequential 0..5
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
